# Evaporation and Demand Plots Update
# Inserts a new row (at row 5) with 0 / 12 / 17 values, and appends two new
# rows of data (rows 9 and 10) at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new row at position 5; existing rows 5-8 shift to 6-9 ---
$ws.Rows(5).Insert()

# --- Step 2: write the new row 5 values (0, 12, 17 repeated across 9 blocks) ---
$row5Values = @("0","12","17","0","12","17","0","12","17","0","12","17","0","12","17","0","12","17","0","12","17","0","12","17","0","12","17")
for ($col = 1; $col -le 27; $col++) {
    $cell = $ws.Cells.Item(5, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $row5Values[$col - 1]
    $cell.Style = "Normal"
}

# --- Step 3: append new row 9 (evaporation figures) ---
$row9Values = @("0","0","0","9142","9142","9142","5485","5485","5485","0","0","0","16237","16237","16237","9742","9742","9742","0","0","0","22331","22331","22331","13398","13398","13398")
for ($col = 1; $col -le 27; $col++) {
    $cell = $ws.Cells.Item(9, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $row9Values[$col - 1]
    $cell.Style = "Normal"
}

# --- Step 4: append new row 10 (demand figures) ---
$row10Values = @("2536404","2536404","2536404","2536404","2536404","2536404","2536404","2536404","2536404","5109197","5109197","5109197","5109197","5109197","5109197","5109197","5109197","5109197","6965058","6965058","6965058","6965058","6965058","6965058","6965058","6965058","6965058")
for ($col = 1; $col -le 27; $col++) {
    $cell = $ws.Cells.Item(10, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $row10Values[$col - 1]
    $cell.Style = "Normal"
}

Write-Host "Edit applied"
